$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2849
$ws.Range("F3").Value = 1153
$ws.Range("F4").Value = 20854
$ws.Range("F6").Value = 2715
$ws.Range("F7").Value = 794
$ws.Range("F9").Value = 500
$ws.Range("F10").Value = 752
$ws.Range("F11").Value = 274
$ws.Range("F14").Value = 105
$ws.Range("F15").Value = 506
$ws.Range("F17").Value = 250
$ws.Range("F18").Value = 11
$ws.Range("F19").Value = 414
$ws.Range("F20").Value = 34
$ws.Range("F23").Value = 16
$ws.Range("F24").Value = 120

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 17
$ws.Range("F3").Value = 28
$ws.Range("F6").Value = 140
$ws.Range("F10").Value = 14
$ws.Range("F14").Value = 136
$ws.Range("F18").Value = 3

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6110
$ws.Range("F3").Value = 691
$ws.Range("F4").Value = 669
$ws.Range("F5").Value = 1521
$ws.Range("F6").Value = 50

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6110
$ws.Range("F3").Value = 691
$ws.Range("F4").Value = 669
$ws.Range("F5").Value = 1521
$ws.Range("F6").Value = 2849
$ws.Range("F7").Value = 1153
$ws.Range("F8").Value = 20855
$ws.Range("F9").Value = 17
$ws.Range("F10").Value = 28
$ws.Range("F14").Value = 2715
$ws.Range("F15").Value = 794
$ws.Range("F16").Value = 140
$ws.Range("F17").Value = 50
$ws.Range("F19").Value = 500
$ws.Range("F20").Value = 752
$ws.Range("F21").Value = 274
$ws.Range("F27").Value = 105
$ws.Range("F28").Value = 14
$ws.Range("F30").Value = 506
$ws.Range("F34").Value = 250
$ws.Range("F35").Value = 136
$ws.Range("F36").Value = 136
$ws.Range("F37").Value = 11
$ws.Range("F38").Value = 414
$ws.Range("F44").Value = 16
$ws.Range("F45").Value = 3
$ws.Range("F50").Value = 120
